$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "85-11=74"
$t.Cell(1,2).Range.Text = "3+74=77"
$t.Cell(1,3).Range.Text = "1+77=78"
$t.Cell(1,4).Range.Text = "87-31=56"
$t.Cell(1,5).Range.Text = "33+47=80"

$t.Cell(2,1).Range.Text = "69-50=19"
$t.Cell(2,2).Range.Text = "36+5=41"
$t.Cell(2,3).Range.Text = "37-28=9"
$t.Cell(2,4).Range.Text = "37+50=87"
$t.Cell(2,5).Range.Text = "60-57=3"

$t.Cell(3,1).Range.Text = "16-9=7"
$t.Cell(3,2).Range.Text = "99-12=87"
$t.Cell(3,3).Range.Text = "19+35=54"
$t.Cell(3,4).Range.Text = "57-26=31"
$t.Cell(3,5).Range.Text = "77-16=61"

$t.Cell(4,1).Range.Text = "8+75=83"
$t.Cell(4,2).Range.Text = "83-7=76"
$t.Cell(4,3).Range.Text = "81-16=65"
$t.Cell(4,4).Range.Text = "23+74=97"
$t.Cell(4,5).Range.Text = "66-34=32"

$t.Cell(5,1).Range.Text = "62-27=35"
$t.Cell(5,2).Range.Text = "30-12=18"
$t.Cell(5,3).Range.Text = "33+50=83"
$t.Cell(5,4).Range.Text = "0+58=58"
$t.Cell(5,5).Range.Text = "40+53=93"

$t.Cell(6,1).Range.Text = "91-41=50"
$t.Cell(6,2).Range.Text = "10+72=82"
$t.Cell(6,3).Range.Text = "91-66=25"
$t.Cell(6,4).Range.Text = "65-47=18"
$t.Cell(6,5).Range.Text = "90-12=78"

$t.Cell(7,1).Range.Text = "50-9=41"
$t.Cell(7,2).Range.Text = "16+71=87"
$t.Cell(7,3).Range.Text = "18+12=30"
$t.Cell(7,4).Range.Text = "14+57=71"
$t.Cell(7,5).Range.Text = "79+15=94"

$t.Cell(8,1).Range.Text = "44+35=79"
$t.Cell(8,2).Range.Text = "92-4=88"
$t.Cell(8,3).Range.Text = "18-7=11"
$t.Cell(8,4).Range.Text = "93-19=74"
$t.Cell(8,5).Range.Text = "6+84=90"

$t.Cell(9,1).Range.Text = "15-11=4"
$t.Cell(9,2).Range.Text = "19+18=37"
$t.Cell(9,3).Range.Text = "88-61=27"
$t.Cell(9,4).Range.Text = "19+42=61"
$t.Cell(9,5).Range.Text = "49+49=98"

$t.Cell(10,1).Range.Text = "97-21=76"
$t.Cell(10,2).Range.Text = "74-6=68"
$t.Cell(10,3).Range.Text = "51+44=95"
$t.Cell(10,4).Range.Text = "67+20=87"
$t.Cell(10,5).Range.Text = "45+15=60"

$t.Cell(11,1).Range.Text = "98-52=46"
$t.Cell(11,2).Range.Text = "99-50=49"
$t.Cell(11,3).Range.Text = "94-63=31"
$t.Cell(11,4).Range.Text = "68+10=78"
$t.Cell(11,5).Range.Text = "46+25=71"

$t.Cell(12,1).Range.Text = "75-61=14"
$t.Cell(12,2).Range.Text = "63-53=10"
$t.Cell(12,3).Range.Text = "21+17=38"
$t.Cell(12,4).Range.Text = "62-7=55"
$t.Cell(12,5).Range.Text = "66-66=0"

$t.Cell(13,1).Range.Text = "74-26=48"
$t.Cell(13,2).Range.Text = "12+67=79"
$t.Cell(13,3).Range.Text = "3+76=79"
$t.Cell(13,4).Range.Text = "52-21=31"
$t.Cell(13,5).Range.Text = "66-7=59"

$t.Cell(14,1).Range.Text = "53-3=50"
$t.Cell(14,2).Range.Text = "34-11=23"
$t.Cell(14,3).Range.Text = "70+10=80"
$t.Cell(14,4).Range.Text = "37+35=72"
$t.Cell(14,5).Range.Text = "98-55=43"

$t.Cell(15,1).Range.Text = "43-29=14"
$t.Cell(15,2).Range.Text = "30+45=75"
$t.Cell(15,3).Range.Text = "72-36=36"
$t.Cell(15,4).Range.Text = "6+84=90"
$t.Cell(15,5).Range.Text = "0+59=59"

$t.Cell(16,1).Range.Text = "20+7=27"
$t.Cell(16,2).Range.Text = "80-5=75"
$t.Cell(16,3).Range.Text = "10+63=73"
$t.Cell(16,4).Range.Text = "78-56=22"
$t.Cell(16,5).Range.Text = "99-20=79"

$t.Cell(17,1).Range.Text = "75-60=15"
$t.Cell(17,2).Range.Text = "35+31=66"
$t.Cell(17,3).Range.Text = "77-2=75"
$t.Cell(17,4).Range.Text = "21+37=58"
$t.Cell(17,5).Range.Text = "22-20=2"

$t.Cell(18,1).Range.Text = "50+13=63"
$t.Cell(18,2).Range.Text = "47-8=39"
$t.Cell(18,3).Range.Text = "55-28=27"
$t.Cell(18,4).Range.Text = "77-32=45"
$t.Cell(18,5).Range.Text = "44-2=42"

$t.Cell(19,1).Range.Text = "10+33=43"
$t.Cell(19,2).Range.Text = "24-19=5"
$t.Cell(19,3).Range.Text = "44+22=66"
$t.Cell(19,4).Range.Text = "70-68=2"
$t.Cell(19,5).Range.Text = "69-55=14"

$t.Cell(20,1).Range.Text = "58+14=72"
$t.Cell(20,2).Range.Text = "69-0=69"
$t.Cell(20,3).Range.Text = "38-11=27"
$t.Cell(20,4).Range.Text = "97-63=34"
$t.Cell(20,5).Range.Text = "26+14=40"
